$wb = $excel.ActiveWorkbook

# Scheduled-runner style data refresh: update hardcoded market/profit values
# per-sheet, per-row, matching the latest pull. Values are literal (no formulas
# in this workbook), so we set .Value directly cell-by-cell.

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1024.1578
$ws.Range("I28").Value = 1036.6111
$ws.Range("J28").Value = 800
$ws.Range("K28").Value = 1036.6111
$ws.Range("L28").Value = 800
$ws.Range("M28").Value = -551.6111000000001
$ws.Range("N28").Value = -1770
# Row 62
$ws.Range("H62").Value = 2700.625
$ws.Range("I62").Value = 3167.5
$ws.Range("K62").Value = 3167.5
$ws.Range("M62").Value = -2543.5
# Row 64
$ws.Range("H64").Value = 4801.8076
$ws.Range("I64").Value = 3456.6843
$ws.Range("K64").Value = 3456.6843
$ws.Range("M64").Value = -3208.6843
# Row 65
$ws.Range("H65").Value = 2700.625
$ws.Range("I65").Value = 3167.5
$ws.Range("K65").Value = 15837.5
$ws.Range("M65").Value = -12717.5
# Row 67
$ws.Range("H67").Value = 4801.8076
$ws.Range("I67").Value = 3456.6843
$ws.Range("K67").Value = 3456.6843
$ws.Range("M67").Value = -2598.6843
# Row 74
$ws.Range("H74").Value = 3452.3333
$ws.Range("I74").Value = 3357
$ws.Range("K74").Value = 3357
$ws.Range("M74").Value = -2421
# Row 77
$ws.Range("H77").Value = 3452.3333
$ws.Range("I77").Value = 3357
$ws.Range("K77").Value = 16785
$ws.Range("M77").Value = -12105
# Row 92
$ws.Range("H92").Value = 1267.1945
$ws.Range("I92").Value = 1138.5862
$ws.Range("J92").Value = 1800
$ws.Range("K92").Value = 1138.5862
$ws.Range("L92").Value = 1800
$ws.Range("M92").Value = 109.4138
$ws.Range("N92").Value = -4296
# Row 98
$ws.Range("H98").Value = 3170.1667
$ws.Range("I98").Value = 3003
$ws.Range("J98").Value = 4006
$ws.Range("K98").Value = 3003
$ws.Range("L98").Value = 4006
$ws.Range("M98").Value = -1505
$ws.Range("N98").Value = -7002
# Row 106
$ws.Range("H106").Value = 336525
$ws.Range("I106").Value = 419407.9
$ws.Range("K106").Value = 419407.9
$ws.Range("M106").Value = -418776.9
# Row 107
$ws.Range("H107").Value = 6188.636
$ws.Range("I107").Value = 7846.875
$ws.Range("J107").Value = 1766.6666
$ws.Range("K107").Value = 7846.875
$ws.Range("L107").Value = 1766.6666
$ws.Range("M107").Value = -5926.875
$ws.Range("N107").Value = -5606.6666
# Row 118
$ws.Range("H118").Value = 3363.6086
$ws.Range("I118").Value = 1344
$ws.Range("J118").Value = 4917.154
$ws.Range("K118").Value = 4032
$ws.Range("L118").Value = 14751.462
$ws.Range("M118").Value = -2375
$ws.Range("N118").Value = -18065.462
# Row 122
$ws.Range("H122").Value = 3170.1667
$ws.Range("I122").Value = 3003
$ws.Range("J122").Value = 4006
$ws.Range("K122").Value = 9009
$ws.Range("L122").Value = 12018
$ws.Range("M122").Value = -6559
$ws.Range("N122").Value = -16918
# Row 138
$ws.Range("H138").Value = 5683906
$ws.Range("J138").Value = 9617773
$ws.Range("L138").Value = 28853319
$ws.Range("N138").Value = -28863599

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10373.822
$ws.Range("I32").Value = 11453.105
$ws.Range("K32").Value = 11453.105
$ws.Range("M32").Value = -11166.105
# Row 63
$ws.Range("H63").Value = 2107
$ws.Range("I63").Value = 2107
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2107
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 2107
$ws.Range("I66").Value = 2107
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10535
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
# Row 122
$ws.Range("H122").Value = 8007.375
$ws.Range("I122").Value = 8098.2
$ws.Range("J122").Value = 7553.25
$ws.Range("K122").Value = 24294.6
$ws.Range("L122").Value = 22659.75
$ws.Range("M122").Value = -21844.6
$ws.Range("N122").Value = -27559.75

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 4700.2964
$ws.Range("I105").Value = 3381.6
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 3381.6
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -1634.6
$ws.Range("N105").Value = -8494

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 2313.3809
$ws.Range("I62").Value = 2342.8125
$ws.Range("J62").Value = 2219.2
$ws.Range("K62").Value = 2342.8125
$ws.Range("L62").Value = 2219.2
$ws.Range("M62").Value = -1718.8125
$ws.Range("N62").Value = -3467.2
# Row 65
$ws.Range("H65").Value = 2313.3809
$ws.Range("I65").Value = 2342.8125
$ws.Range("J65").Value = 2219.2
$ws.Range("K65").Value = 11714.0625
$ws.Range("L65").Value = 11096
$ws.Range("M65").Value = -8594.0625
$ws.Range("N65").Value = -17336

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 530
$ws.Range("I5").Value = 487.27274
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1461.81822
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1349.81822
$ws.Range("N5").Value = -3224
# Row 68
$ws.Range("H68").Value = 1369.4286
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1369.4286
$ws.Range("K68").Value = 0
$ws.Range("N68").Value = -5730.2858
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 1369.4286
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1369.4286
$ws.Range("K71").Value = 0
$ws.Range("N71").Value = -20436.8574
$ws.Range("M71").ClearContents()
# Row 131
$ws.Range("H131").Value = 960.913
$ws.Range("I131").Value = 627
$ws.Range("J131").Value = 1107
$ws.Range("K131").Value = 1881
$ws.Range("L131").Value = 3321
$ws.Range("M131").Value = 3159
$ws.Range("N131").Value = -13401
# Row 135
$ws.Range("H135").Value = 530
$ws.Range("I135").Value = 487.27274
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 4385.45466
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -1850.45466
$ws.Range("N135").Value = -14070
# Row 137
$ws.Range("H137").Value = 6177870.5
$ws.Range("I137").Value = 18522424
$ws.Range("J137").Value = 5594.3887
$ws.Range("K137").Value = 55567272
$ws.Range("L137").Value = 16783.1661
$ws.Range("M137").Value = -55562172
$ws.Range("N137").Value = -26983.1661

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4463.2
$ws.Range("I7").Value = 4905
$ws.Range("K7").Value = 4905
$ws.Range("M7").Value = -4793
# Row 16
$ws.Range("H16").Value = 2364.96
$ws.Range("I16").Value = 2357.4783
$ws.Range("J16").Value = 2451
$ws.Range("K16").Value = 2357.4783
$ws.Range("L16").Value = 2451
$ws.Range("M16").Value = -2187.4783
$ws.Range("N16").Value = -2791
# Row 46
$ws.Range("H46").Value = 1486.6666
$ws.Range("I46").Value = 1375
$ws.Range("J46").Value = 1614.2858
$ws.Range("K46").Value = 1375
$ws.Range("L46").Value = 1614.2858
$ws.Range("M46").Value = -1187
$ws.Range("N46").Value = -1990.2858
# Row 93
$ws.Range("H93").Value = 1095.3214
$ws.Range("I93").Value = 840.0526
$ws.Range("J93").Value = 1634.2222
$ws.Range("K93").Value = 840.0526
$ws.Range("L93").Value = 1634.2222
$ws.Range("M93").Value = 407.9474
$ws.Range("N93").Value = -4130.2222
# Row 126
$ws.Range("H126").Value = 4463.2
$ws.Range("I126").Value = 4905
$ws.Range("K126").Value = 14715
$ws.Range("M126").Value = -12245

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1906.4667
$ws.Range("I100").Value = 2459.6667
$ws.Range("J100").Value = 1076.6666
$ws.Range("K100").Value = 4919.3334
$ws.Range("L100").Value = 2153.3332
$ws.Range("M100").Value = -4378.3334
$ws.Range("N100").Value = -3235.3332
# Row 122
$ws.Range("H122").Value = 2649.0417
$ws.Range("I122").Value = 2555.4707
$ws.Range("J122").Value = 2876.2856
$ws.Range("K122").Value = 7666.4121
$ws.Range("L122").Value = 8628.856800000001
$ws.Range("M122").Value = -5216.4121
$ws.Range("N122").Value = -13528.8568
